# RPAR_holdings.xlsx update:
#  - Bump the "as of" date in the confidentiality footnote from 2021-05-18 to 2021-05-19
#  - Refresh the Weight / Percent Change figures for every holding (rows 2-15)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected; unprotect so the cells can be written, then
# restore protection once the data has been refreshed.
$ws.Unprotect()

# --- Footnote date bump -----------------------------------------------
$footnote = $ws.Range("A18").Value()
$ws.Range("A18").Value = $footnote -replace "2021-05-18", "2021-05-19"
# Re-fit the row height (auto-wrap text) back to its natural size so we
# don't leave behind an explicit row-height override from the text change.
$ws.Rows(18).AutoFit()

# --- Weight (col D) / Percent Change (col E) refresh --------------------
$ws.Range("D2").Value = 0.0561181320866728
$ws.Range("E2").Value = -0.002905342080599782

$ws.Range("D3").Value = 0.02375093437084382
$ws.Range("E3").Value = -0.008913001356326222

$ws.Range("D4").Value = 0.03111085106145185
$ws.Range("E4").Value = -0.002662609357169998

$ws.Range("D5").Value = 0.0333263865011436
$ws.Range("E5").Value = -0.02548725637181415

$ws.Range("D6").Value = 0.03852270463625445
$ws.Range("E6").Value = -0.03297153389309859

$ws.Range("D7").Value = 0.0191811459906369
$ws.Range("E7").Value = -0.02327255278310925

$ws.Range("D8").Value = 0.004299585096612467
$ws.Range("E8").Value = 0.007798165137614443

$ws.Range("D9").Value = 0.006803074400572753
$ws.Range("E9").Value = -0.00676459219172787

$ws.Range("D10").Value = 0.07325072958173716
$ws.Range("E10").Value = 0.001077005923532548

$ws.Range("D11").Value = 0.07332962105139977
$ws.Range("E11").Value = 0.001075847229693228

$ws.Range("D12").Value = 0.1435193616102275
$ws.Range("E12").Value = -0.002418645558487142

$ws.Range("D13").Value = 0.3823632860137916
$ws.Range("E13").Value = -0.002981410031567933

$ws.Range("D14").Value = 0.1144241875986554
$ws.Range("E14").Value = -0.01206563706563701

$ws.Range("D15").Value = 1
$ws.Range("E15").Value = -0.005745928707093673

# --- Restore sheet protection (matches original locked-sheet state) ----
$ws.Protect("lido", $true, $true, $true)
